$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4079.8
$ws.Range("J51").Value = 4107.6924
$ws.Range("L51").Value = 4107.6924
$ws.Range("N51").Value = -5075.6924
$ws.Range("H80").Value = 8330.053
$ws.Range("I80").Value = 995.2857
$ws.Range("J80").Value = 12608.667
$ws.Range("K80").Value = 2985.8571
$ws.Range("L80").Value = 37826.001
$ws.Range("M80").Value = -1987.8571
$ws.Range("N80").Value = -39822.001
$ws.Range("H83").Value = 8330.053
$ws.Range("I83").Value = 995.2857
$ws.Range("J83").Value = 12608.667
$ws.Range("K83").Value = 8957.5713
$ws.Range("L83").Value = 113478.003
$ws.Range("M83").Value = -3965.5713
$ws.Range("N83").Value = -123462.003
$ws.Range("H116").Value = 6215.5386
$ws.Range("I116").Value = 4003
$ws.Range("J116").Value = 6879.3
$ws.Range("K116").Value = 4003
$ws.Range("L116").Value = 6879.3
$ws.Range("M116").Value = -561
$ws.Range("N116").Value = -13763.3
$ws.Range("H137").Value = 551399
$ws.Range("I137").Value = 1827.7333
$ws.Range("J137").Value = 1300814.4
$ws.Range("K137").Value = 5483.199900000001
$ws.Range("L137").Value = 3902443.2
$ws.Range("M137").Value = -2933.199900000001
$ws.Range("N137").Value = -3907543.2
$ws.Range("H138").Value = 2353.8293
$ws.Range("I138").Value = 1683.2307
$ws.Range("J138").Value = 2665.1785
$ws.Range("K138").Value = 5049.6921
$ws.Range("L138").Value = 7995.5355
$ws.Range("M138").Value = 90.30789999999979
$ws.Range("N138").Value = -18275.5355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6778.48
$ws.Range("J32").Value = 27543.215
$ws.Range("L32").Value = 27543.215
$ws.Range("N32").Value = -28117.215
$ws.Range("H132").Value = 2466.4285
$ws.Range("I132").Value = 2025.9333
$ws.Range("K132").Value = 6077.7999
$ws.Range("M132").Value = -3547.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2596.0833
$ws.Range("I86").Value = 1252
$ws.Range("K86").Value = 1252
$ws.Range("M86").Value = -129
$ws.Range("H89").Value = 2596.0833
$ws.Range("I89").Value = 1252
$ws.Range("K89").Value = 6260
$ws.Range("M89").Value = -644
$ws.Range("H130").Value = 79780
$ws.Range("J130").Value = 79780
$ws.Range("L130").Value = 79780
$ws.Range("N130").Value = -89820
$ws.Range("H135").Value = 105561.43
$ws.Range("J135").Value = 105561.43
$ws.Range("L135").Value = 105561.43
$ws.Range("N135").Value = -115701.43
$ws.Range("H138").Value = 71989.75
$ws.Range("J138").Value = 71989.75
$ws.Range("L138").Value = 71989.75
$ws.Range("N138").Value = -82269.75
$ws.Range("H139").Value = 49999
$ws.Range("J139").Value = 49999
$ws.Range("L139").Value = 49999
$ws.Range("N139").Value = -60279
$ws.Range("H140").Value = 55969.082
$ws.Range("J140").Value = 55969.082
$ws.Range("L140").Value = 55969.082
$ws.Range("N140").Value = -66329.08199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4058.491
$ws.Range("I31").Value = 2109.5806
$ws.Range("J31").Value = 6575.8335
$ws.Range("K31").Value = 2109.5806
$ws.Range("L31").Value = 6575.8335
$ws.Range("M31").Value = -1814.5806
$ws.Range("N31").Value = -7165.8335
$ws.Range("H34").Value = 4058.491
$ws.Range("I34").Value = 2109.5806
$ws.Range("J34").Value = 6575.8335
$ws.Range("K34").Value = 2109.5806
$ws.Range("L34").Value = 6575.8335
$ws.Range("M34").Value = -1907.5806
$ws.Range("N34").Value = -6979.8335
$ws.Range("H43").Value = 32121
$ws.Range("J43").Value = 32121
$ws.Range("L43").Value = 32121
$ws.Range("N43").Value = -32489
$ws.Range("H52").Value = 84995.664
$ws.Range("I52").Value = 54987
$ws.Range("J52").Value = 100000
$ws.Range("K52").Value = 54987
$ws.Range("L52").Value = 100000
$ws.Range("M52").Value = -54693
$ws.Range("N52").Value = -100588
$ws.Range("H86").Value = 125003864
$ws.Range("I86").Value = 142861180
$ws.Range("K86").Value = 142861180
$ws.Range("M86").Value = -142860057
$ws.Range("H89").Value = 125003864
$ws.Range("I89").Value = 142861180
$ws.Range("K89").Value = 714305900
$ws.Range("M89").Value = -714300284
$ws.Range("H101").Value = 32121
$ws.Range("J101").Value = 32121
$ws.Range("L101").Value = 32121
$ws.Range("N101").Value = -38611
$ws.Range("H122").Value = 2668.4546
$ws.Range("I122").Value = 2142.5
$ws.Range("K122").Value = 6427.5
$ws.Range("M122").Value = -3977.5
$ws.Range("H132").Value = 2162.3333
$ws.Range("I132").Value = 1617.1538
$ws.Range("J132").Value = 3579.8
$ws.Range("K132").Value = 4851.4614
$ws.Range("L132").Value = 10739.4
$ws.Range("M132").Value = -2321.4614
$ws.Range("N132").Value = -15799.4
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 171815.73
$ws.Range("J141").Value = 171815.73
$ws.Range("L141").Value = 171815.73
$ws.Range("N141").Value = -182175.73

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 52756.145
$ws.Range("J37").Value = 52756.145
$ws.Range("L37").Value = 158268.435
$ws.Range("N37").Value = -158492.435
$ws.Range("H131").Value = 480191.25
$ws.Range("J131").Value = 668960.9
$ws.Range("L131").Value = 2006882.7
$ws.Range("N131").Value = -2016962.7
$ws.Range("H133").Value = 7034.8335
$ws.Range("J133").Value = 8074.0835
$ws.Range("L133").Value = 24222.2505
$ws.Range("N133").Value = -34342.25049999999
$ws.Range("H134").Value = 1290
$ws.Range("I134").Value = 1290
$ws.Range("K134").Value = 3870
$ws.Range("M134").Value = 1200
$ws.Range("H136").Value = 1377.8334
$ws.Range("I136").Value = 1246.8
$ws.Range("K136").Value = 3740.4
$ws.Range("M136").Value = 1359.6
$ws.Range("H137").Value = 3040.375
$ws.Range("I137").Value = 1278
$ws.Range("K137").Value = 3834
$ws.Range("M137").Value = 1266
$ws.Range("H139").Value = 12506651
$ws.Range("I139").Value = 19232930
$ws.Range("J139").Value = 14991.429
$ws.Range("K139").Value = 57698790
$ws.Range("L139").Value = 44974.287
$ws.Range("M139").Value = -57693650
$ws.Range("N139").Value = -55254.287
$ws.Range("H141").Value = 45458084
$ws.Range("I141").Value = 50002892
$ws.Range("K141").Value = 150008676
$ws.Range("M141").Value = -150003496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 41957.035
$ws.Range("I70").Value = 55942.85
$ws.Range("J70").Value = 6992.5
$ws.Range("K70").Value = 55942.85
$ws.Range("L70").Value = 6992.5
$ws.Range("M70").Value = -55672.85
$ws.Range("N70").Value = -7532.5
$ws.Range("H73").Value = 41957.035
$ws.Range("I73").Value = 55942.85
$ws.Range("J73").Value = 6992.5
$ws.Range("K73").Value = 55942.85
$ws.Range("L73").Value = 6992.5
$ws.Range("M73").Value = -55006.85
$ws.Range("N73").Value = -8864.5
$ws.Range("H122").Value = 11152.259
$ws.Range("I122").Value = 3396
$ws.Range("J122").Value = 45279.8
$ws.Range("K122").Value = 10188
$ws.Range("L122").Value = 135839.4
$ws.Range("M122").Value = -7738
$ws.Range("N122").Value = -140739.4
$ws.Range("H132").Value = 1862.0244
$ws.Range("I132").Value = 1740.9697
$ws.Range("J132").Value = 2361.375
$ws.Range("K132").Value = 5222.909100000001
$ws.Range("L132").Value = 7084.125
$ws.Range("M132").Value = -2692.909100000001
$ws.Range("N132").Value = -12144.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2398.8
$ws.Range("I16").Value = 2398.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2398.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2228.8
$ws.Range("N16").Value = ""
$ws.Range("H40").Value = 2926936
$ws.Range("J40").Value = 7939607
$ws.Range("L40").Value = 7939607
$ws.Range("N40").Value = -7939879
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -798
$ws.Range("H97").Value = 13267.25
$ws.Range("J97").Value = 13267.25
$ws.Range("L97").Value = 13267.25
$ws.Range("N97").Value = -15249.25
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H136").Value = 6340.727
$ws.Range("I136").Value = 6366.091
$ws.Range("J136").Value = 6315.364
$ws.Range("K136").Value = 19098.273
$ws.Range("L136").Value = 18946.092
$ws.Range("M136").Value = -16548.273
$ws.Range("N136").Value = -24046.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 37991
$ws.Range("J59").Value = 37991
$ws.Range("L59").Value = 37991
$ws.Range("N59").Value = -39467
$ws.Range("H101").Value = 61900
$ws.Range("J101").Value = 61900
$ws.Range("L101").Value = 61900
$ws.Range("N101").Value = -68390
$ws.Range("H103").Value = 28666.666
$ws.Range("J103").Value = 28666.666
$ws.Range("L103").Value = 28666.666
$ws.Range("N103").Value = -31010.666
$ws.Range("H122").Value = 4213
$ws.Range("I122").Value = 2678.6
$ws.Range("J122").Value = 5491.6665
$ws.Range("K122").Value = 8035.799999999999
$ws.Range("L122").Value = 16474.9995
$ws.Range("M122").Value = -5585.799999999999
$ws.Range("N122").Value = -21374.9995
$ws.Range("H132").Value = 1674455.1
$ws.Range("I132").Value = 2097.3684
$ws.Range("K132").Value = 6292.1052
$ws.Range("M132").Value = -3762.1052
